$wb = $excel.ActiveWorkbook

# --- New sheet "HelloWorld" appended after the last existing sheet,
#     carrying the old Hello/World content that used to live on Summary ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$helloSheet = $wb.Worksheets.Add($null, $lastSheet)
$helloSheet.Name = "HelloWorld"
$helloSheet.Range("A1").Value = "Hello"
$helloSheet.Range("B1").Value = "World"

# --- Summary sheet: replace Hello/World text with numbers + SUM formula,
#     add an empty row 2 with A2:B2 merged ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("A1").Value = 10
$summary.Range("B1").Value = 20
$summary.Range("C1").Formula = "=SUM(A1,B1)"
$summary.Range("A2:B2").Merge()
